$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") should look exactly like the
# existing header cells (bold, centered, thin box border) - copy the
# formatting from H1 (an existing header cell) rather than rebuilding it
# property by property, then stamp in the text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new I (I0) and J (IF) columns, rows 2-25.
$iValues = @(5, 6, 8, 5, 8, 8, 6, 5, 7, 8, 9, 9, 7, 8, 5, 5, 2, 7, 6, 7, 8, 3, 4, 4)
$jValues = @(7, 6, 8, 7, 8, 8, 7, 6, 7, 8, 9, 9, 7, 9, 6, 5, 2, 8, 6, 7, 8, 3, 4, 4)

for ($r = 2; $r -le 25; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
